$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C53").Value = "thumbs_up"
$ws.Range("D53").Value = "TINYINT"
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "NOT NULL"
$ws.Range("H53").Value = "是否喜欢"
$ws.Range("I53").Value = "类似为点赞，默认值为0"

$ws.Range("I53").Select()
